$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.508.81"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.676.86"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "219.91"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "0.5316"
$ws.Range("E6").Value = "  +1.93%  "
$ws.Range("D8").Value = "0.2693"
$ws.Range("E8").Value = "  +3.28%  "
$ws.Range("D9").Value = "0.06400"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").Value = "21.74"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").Value = "0.07796"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "1.696.72"
$ws.Range("E12").Value = "  +2.85%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").Value = "0.5581"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "0.0₅8352"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "26.531.22"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D19").Value = "4.788"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "192.58"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "0.1276"
$ws.Range("E24").Value = "  +5.31%  "
$ws.Range("D25").Value = "139.29"
$ws.Range("E25").Value = "  -4.87%  "
$ws.Range("D26").Value = "7.426"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").Value = "1.436"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("E29").Value = "  +7.17%  "
$ws.Range("D30").Value = "1.289"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").Value = "3.607"
$ws.Range("E31").Value = "  +5.51%  "
$ws.Range("D32").Value = "3.440"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "1.694"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("D35").Value = "0.6154"
$ws.Range("E35").Value = "  +8.63%  "
$ws.Range("D36").Value = "2.426"
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").Value = "2.781"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "0.01630"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").Value = "6.097"
$ws.Range("E39").Value = "  +4.52%  "
$ws.Range("D40").Value = "1.094.53"
$ws.Range("E40").Value = "  +6.23%  "
$ws.Range("D41").Value = "0.8619"
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").Value = "100.62"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "1.823.65"
$ws.Range("E44").Value = "  +1.46%  "
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("D46").Value = "58.67"
$ws.Range("E46").Value = "  +4.83%  "
$ws.Range("D47").Value = "8.188"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "0.9979"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").Value = "1.513"
$ws.Range("E49").Value = "  +9.39%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "6.024"
$ws.Range("E51").Value = "  +1.62%  "
